# kpinhanvien.xlsx: add a new "ĐIỂM TRỪ (%)" column (D) to the KPI table,
# matching the formatting already used by column C, and fill in the one
# data point the author entered (row 3 -> 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- New column D: copy column C's per-cell formatting (borders, fill,
#     alignment) row by row, so the new header/data cells look identical
#     to their column-C neighbours. ---
for ($r = 1; $r -le 8; $r++) {
    $src = $ws.Cells.Item($r, 3)
    $dst = $ws.Cells.Item($r, 4)
    $src.Copy()
    $dst.PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = $false

# --- Header + the single populated data value in the new column ---
$ws.Range("D1").Value = "ĐIỂM TRỪ (%)"
$ws.Range("D3").Value = 5

# --- Column width for D, matching column C (20 chars wide) ---
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# --- Cursor/selection ends up on E3 after the edit ---
$ws.Range("E3").Select()

# --- Best-effort: the author's window also got resized/moved on their
#     desktop (workbookView + absPath in the XML). There is no headless
#     equivalent of a real on-screen window in this runtime, so these
#     calls are harmless no-ops here, but are included for parity with
#     what the author did interactively in Excel. ---
$excel.ActiveWindow.Width = 27405
$excel.ActiveWindow.Height = 12915
